# Actualización automática 2025-12-05 14:30:07
# Insert a new client "CARRION CARRION STEPHANIE DAYANA" (with all-zero
# sales figures) right before "CHASIQUIZA CAMPAÑA JOSE LUIS" (alphabetical
# order) on both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. This
# pushes every following data row down by one, and the trailing
# summary/total row (previously row 23, now row 24) needs its "X de 21"
# labels bumped to "X de 22" to reflect the new client count.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "VENTAS POR GRUPO" (columns A:R, data rows 2-23, totals row 24) ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(6).Insert()
$ws1.Range("A6").Value = "HIDALGO HIDALGO PEDRO GUSTAVO"
$ws1.Range("B6").Value = "CARRION CARRION STEPHANIE DAYANA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(6, $col).Value = 0
}

# Bump the "X de 21" -> "X de 22" labels on the (now shifted) totals row.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(24, $col)
    $cell.Value = ($cell.Value2 -replace "de 21", "de 22")
}

# ---- Sheet 2: "VENTA MENSUAL" (columns A:G, data rows 2-23, totals row 24) ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(6).Insert()
$ws2.Range("A6").Value = "HIDALGO HIDALGO PEDRO GUSTAVO"
$ws2.Range("B6").Value = "CARRION CARRION STEPHANIE DAYANA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(6, $col).Value = 0
}
